$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 2

$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2.4
$ws.Range("T3").Value = 1.53
$ws.Range("U3").Value = 4
$ws.Range("V3").Value = 1.25
$ws.Range("W3").Value = 4.5
$ws.Range("X3").Value = 1.18

$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AF4").Value = 10
$ws.Range("AI4").Value = 7.5
$ws.Range("AL4").Value = 81
$ws.Range("AN4").Value = 13
$ws.Range("AO4").Value = 34
$ws.Range("AQ4").Value = 81
$ws.Range("AR4").Value = 51
$ws.Range("G4").Value = 1.53
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 2.2
$ws.Range("L4").Value = 7

$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 3

$ws.Range("AA6").Value = 1.73
$ws.Range("AB6").Value = 2
$ws.Range("AC6").Value = 10
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 34
$ws.Range("AI6").Value = 10
$ws.Range("AM6").Value = 201
$ws.Range("AN6").Value = 8
$ws.Range("G6").Value = 3.3
$ws.Range("H6").Value = 3.3
$ws.Range("J6").Value = 3.75
$ws.Range("K6").Value = 2.1
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 3.4
$ws.Range("S6").Value = 2.05
$ws.Range("T6").Value = 1.8
$ws.Range("W6").Value = 3.5
$ws.Range("X6").Value = 1.29
$ws.Range("Y6").Value = 1.4
$ws.Range("Z6").Value = 2.75

$ws.Range("AC7").Value = 10
$ws.Range("AD7").Value = 17
$ws.Range("AI7").Value = 10
$ws.Range("AN7").Value = 8
$ws.Range("AP7").Value = 9
$ws.Range("G7").Value = 3.25
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 2.2
$ws.Range("L7").Value = 2.88
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 11
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 3.5
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 1.85
$ws.Range("Y7").Value = 1.4
$ws.Range("Z7").Value = 2.75

$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 18
$ws.Range("AE8").Value = 11.75
$ws.Range("AF8").Value = 45
$ws.Range("AG8").Value = 32
$ws.Range("AN8").Value = 7.4
$ws.Range("AO8").Value = 9.75
$ws.Range("AQ8").Value = 18
$ws.Range("AR8").Value = 16
$ws.Range("G8").Value = 3.45
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 2.02
$ws.Range("J8").Value = 3.95
$ws.Range("L8").Value = 2.6
$ws.Range("M8").Value = 1.06
$ws.Range("Z8").Value = 2.75

$ws.Range("AA9").Value = 1.62
$ws.Range("AB9").Value = 2.15
$ws.Range("AE9").Value = 29
$ws.Range("AF9").Value = 250
$ws.Range("AG9").Value = 90
$ws.Range("AH9").Value = 60
$ws.Range("AI9").Value = 30
$ws.Range("AJ9").Value = 13.5
$ws.Range("AK9").Value = 19
$ws.Range("AL9").Value = 55
$ws.Range("AM9").Value = 300
$ws.Range("AN9").Value = 12.5
$ws.Range("AO9").Value = 9
$ws.Range("AP9").Value = 9.5
$ws.Range("AQ9").Value = 9.25
$ws.Range("AR9").Value = 9.75
$ws.Range("AS9").Value = 19.5
$ws.Range("G9").Value = 9.25
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 1.24
$ws.Range("J9").Value = 7
$ws.Range("K9").Value = 2.9
$ws.Range("L9").Value = 1.6
$ws.Range("P9").Value = 6.2
$ws.Range("S9").Value = 1.3
$ws.Range("T9").Value = 3.25
$ws.Range("W9").Value = 1.72
$ws.Range("X9").Value = 2
